# Added a poll and amended two old polls
# (column H = "CAN" national numbers on Sheet1's poll table)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New poll added: party vote shares for CAN (H4:H7) change value,
#     and H8:H9 keep their values but get re-entered (style flips from
#     the explicitly-font-applied style 1 to the default style 0, same
#     as what happens when a cell is retyped in the regional/default
#     font rather than inheriting the sheet's original explicit style).
$ws.Range("H4").Value = 35
$ws.Range("H4").Font.Name = "Arial"

$ws.Range("H5").Value = 30
$ws.Range("H5").Font.Name = "Arial"

$ws.Range("H6").Value = 20
$ws.Range("H6").Font.Name = "Arial"

$ws.Range("H7").Value = 7
$ws.Range("H7").Font.Name = "Arial"

$ws.Range("H8").Value = 6
$ws.Range("H8").Font.Name = "Arial"

$ws.Range("H9").Value = 2
$ws.Range("H9").Font.Name = "Arial"

# --- Two old polls amended: sample-size totals (nw / nu row) updated.
$ws.Range("H10").Value = 1238
$ws.Range("H11").Value = 1242

# --- Cursor left on I4 (next empty column) after the edits.
$ws.Range("I4").Select() | Out-Null
